# Actualización 10 de Mayo
# Updates the "Estadisticos 2P" and "Estadisticos Final" sheets with the
# latest grade statistics.

$wb = $excel.ActiveWorkbook

# --- Sheet "Estadisticos 2P" ---------------------------------------------
$ws2P = $wb.Worksheets.Item("Estadisticos 2P")

# Row 2
$ws2P.Range("D2").Value = 0
$ws2P.Range("E2").Value = 0
$ws2P.Range("F2").Value = 39
$ws2P.Range("G2").Value = 100
$ws2P.Range("H2").Value = 9.199999999999999

# Row 3
$ws2P.Range("D3").Value = 0
$ws2P.Range("E3").Value = 0
$ws2P.Range("F3").Value = 39
$ws2P.Range("G3").Value = 100
$ws2P.Range("H3").Value = 9.199999999999999

# Row 4
$ws2P.Range("D4").Value = 0
$ws2P.Range("E4").Value = 0
$ws2P.Range("F4").Value = 37
$ws2P.Range("G4").Value = 100
$ws2P.Range("H4").Value = 8.6

# --- Sheet "Estadisticos Final" ------------------------------------------
$wsFinal = $wb.Worksheets.Item("Estadisticos Final")

$wsFinal.Range("H2").Value = 8.9
$wsFinal.Range("H3").Value = 8.9
$wsFinal.Range("H4").Value = 8.699999999999999
